$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.319.74"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.677.14"
$ws.Range("E3").Value = "  +3.85%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'612.79"
$ws.Range("E5").Value = "  +4.67%  "
$ws.Range("D6").Value = "'143.44"
$ws.Range("E6").Value = "  -0.82%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.588"
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("D9").Value = "2.677.28"
$ws.Range("E9").Value = "  +3.94%  "
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").Value = "'5.61"
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "'0.152"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("E13").Value = "  +3.26%  "
$ws.Range("D14").Value = "'27.40"
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("D15").Value = "3.158.91"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "63.244.52"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.670.92"
$ws.Range("E18").Value = "  +3.82%  "
$ws.Range("E19").Value = "  +3.45%  "
$ws.Range("D20").Value = "'342.53"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("E21").Value = "  +1.80%  "
$ws.Range("E22").Value = "  +3.33%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").Value = "'67.28"
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").Value = "'1.64"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D26").Value = "'1.54"
$ws.Range("E26").Value = "  -4.66%  "
$ws.Range("D27").Value = "'8.67"
$ws.Range("E27").Value = "  +4.83%  "
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("D29").Value = "'542.33"
$ws.Range("E29").Value = "  +15.92%  "
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("E31").Value = "  -1.13%  "
$ws.Range("D32").Value = "'2.07"
$ws.Range("E32").Value = "  +7.19%  "
$ws.Range("E33").Value = "  +7.64%  "
$ws.Range("D34").Value = "0.0₃0807"
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").Value = "'172.13"
$ws.Range("E35").Value = "  -2.50%  "
$ws.Range("D36").Value = "'5.17"
$ws.Range("E36").Value = "  +12.87%  "
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.08%  "
$ws.Range("D39").Value = "'19.23"
$ws.Range("E39").Value = "  +2.00%  "
$ws.Range("E40").Value = "  +10.17%  "
$ws.Range("D41").Value = "'176.60"
$ws.Range("E41").Value = "  +11.38%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").Value = "'22.28"
$ws.Range("E44").Value = "  +3.96%  "
$ws.Range("E45").Value = "  +6.40%  "
$ws.Range("D46").Value = "'0.636"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'0.0965"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'18.86"
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("D51").Value = "'11.32"
$ws.Range("E51").Value = "  -0.68%  "
